# ji/update debbuging errors codes
# Rename the "PR_B" sub-model tag from PR_B_2 -> PR_B_Y2 in the
# name_process_model column (B) for the PR_B / rows 7-11 block, and
# update the active-cell selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7:B11").Value = "PR_B_Y2"

$ws.Range("H4").Select()
